$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$block1 = New-Object 'object[,]' 24,5
$block1[0,0] = 1.02
$block1[0,1] = 1.029546641954758
$block1[0,2] = 1.037859051853081
$block1[0,3] = 1.038796045446549
$block1[0,4] = 1.048414767778013
$block1[1,0] = 1.02
$block1[1,1] = 1.030355857651603
$block1[1,2] = 1.038354697915657
$block1[1,3] = 1.039525645344808
$block1[1,4] = 1.0492504195831
$block1[2,0] = 1.02
$block1[2,1] = 1.030880201427476
$block1[2,2] = 1.03867572982534
$block1[2,3] = 1.039998796304486
$block1[2,4] = 1.049792352304787
$block1[3,0] = 1.02
$block1[3,1] = 1.031100807870976
$block1[3,2] = 1.03881076523843
$block1[3,3] = 1.040197958661639
$block1[3,4] = 1.050020468462421
$block1[4,0] = 1.02
$block1[4,1] = 1.031137858717636
$block1[4,2] = 1.038833442536147
$block1[4,3] = 1.040231413483867
$block1[4,4] = 1.050058786976443
$block1[5,0] = 1.02
$block1[5,1] = 1.03088314850666
$block1[5,2] = 1.038677533888469
$block1[5,3] = 1.040001456542947
$block1[5,4] = 1.049795399276409
$block1[6,0] = 1.02
$block1[6,1] = 1.029819968581593
$block1[6,2] = 1.038026491180804
$block1[6,3] = 1.039042398365352
$block1[6,4] = 1.048696928411944
$block1[7,0] = 1.02
$block1[7,1] = 1.027952156081785
$block1[7,2] = 1.036881777273
$block1[7,3] = 1.037360550257118
$block1[7,4] = 1.046770641905514
$block1[8,0] = 1.02
$block1[8,1] = 1.026710848709719
$block1[8,2] = 1.036120439239493
$block1[8,3] = 1.036244902902232
$block1[8,4] = 1.045492867245414
$block1[9,0] = 1.02
$block1[9,1] = 1.026174296263958
$block1[9,2] = 1.035791224823524
$block1[9,3] = 1.035763164324941
$block1[9,4] = 1.044941124753881
$block1[10,0] = 1.02
$block1[10,1] = 1.025975140079289
$block1[10,2] = 1.035669009603099
$block1[10,3] = 1.035584429092806
$block1[10,4] = 1.044736416880705
$block1[11,0] = 1.02
$block1[11,1] = 1.0260178532444
$block1[11,2] = 1.035695221996813
$block1[11,3] = 1.035622759146855
$block1[11,4] = 1.044780316797775
$block1[12,0] = 1.02
$block1[12,1] = 1.026157831018337
$block1[12,2] = 1.035781121042057
$block1[12,3] = 1.035748385843384
$block1[12,4] = 1.044924198754341
$block1[13,0] = 1.02
$block1[13,1] = 1.026244094984583
$block1[13,2] = 1.035834055572804
$block1[13,3] = 1.035825815678857
$block1[13,4] = 1.045012880228216
$block1[14,0] = 1.02
$block1[14,1] = 1.026746477924031
$block1[14,2] = 1.03614229777321
$block1[14,3] = 1.03627690279342
$block1[14,4] = 1.045529517255302
$block1[15,0] = 1.02
$block1[15,1] = 1.027061862883689
$block1[15,2] = 1.036335771895383
$block1[15,3] = 1.03656021907952
$block1[15,4] = 1.045854004453309
$block1[16,0] = 1.02
$block1[16,1] = 1.027245912348691
$block1[16,2] = 1.036448665355411
$block1[16,3] = 1.036725602227611
$block1[16,4] = 1.046043421019286
$block1[17,0] = 1.02
$block1[17,1] = 1.027308683799229
$block1[17,2] = 1.03648716638552
$block1[17,3] = 1.036782015520865
$block1[17,4] = 1.046108032326051
$block1[18,0] = 1.02
$block1[18,1] = 1.027028015662329
$block1[18,2] = 1.036315009461744
$block1[18,3] = 1.036529808508801
$block1[18,4] = 1.045819174654539
$block1[19,0] = 1.02
$block1[19,1] = 1.0261166070831
$block1[19,2] = 1.035755823974349
$block1[19,3] = 1.035711386267116
$block1[19,4] = 1.044881822631653
$block1[20,0] = 1.02
$block1[20,1] = 1.025544398215874
$block1[20,2] = 1.035404645983492
$block1[20,3] = 1.03519799255126
$block1[20,4] = 1.044293826276242
$block1[21,0] = 1.02
$block1[21,1] = 1.025847657703597
$block1[21,2] = 1.035590773086371
$block1[21,3] = 1.035470039680425
$block1[21,4] = 1.044605405208278
$block1[22,0] = 1.02
$block1[22,1] = 1.027043309484537
$block1[22,2] = 1.036324390978358
$block1[22,3] = 1.036543549336674
$block1[22,4] = 1.045834912282278
$block1[23,0] = 1.02
$block1[23,1] = 1.028434350488041
$block1[23,2] = 1.037177403902733
$block1[23,3] = 1.037794372395094
$block1[23,4] = 1.047267511261249
$ws.Range("B2:F25").Value = $block1

$block2 = New-Object 'object[,]' 24,6
$block2[0,0] = 1.03706768453043
$block2[0,1] = 1.034692816641452
$block2[0,2] = 1.040648735333751
$block2[0,3] = 1.041583059497417
$block2[0,4] = 1.051174678021466
$block2[0,5] = 1.015562995983377
$block2[1,0] = 1.037212946084829
$block2[1,1] = 1.035143827496455
$block2[1,2] = 1.040955107276626
$block2[1,3] = 1.042122958772917
$block2[1,4] = 1.051822309373589
$block2[1,5] = 1.015713873211963
$block2[2,0] = 1.037306026020034
$block2[2,1] = 1.035435721874545
$block2[2,2] = 1.041152946656341
$block2[2,3] = 1.042472690340971
$block2[2,4] = 1.052241930750347
$block2[2,5] = 1.015811481069486
$block2[3,0] = 1.037344937230819
$block2[3,1] = 1.035558447391879
$block2[3,2] = 1.041236020296189
$block2[3,3] = 1.042619807120453
$block2[3,4] = 1.052418472048691
$block2[3,5] = 1.015852510197092
$block2[4,0] = 1.037351457694583
$block2[4,1] = 1.035579054278153
$block2[4,2] = 1.041249962936965
$block2[4,3] = 1.042644513887404
$block2[4,4] = 1.052448121827289
$block2[4,5] = 1.015859398848458
$block2[5,0] = 1.037306546817364
$block2[5,1] = 1.035437361688977
$block2[5,2] = 1.041154057077326
$block2[5,3] = 1.04247465577228
$block2[5,4] = 1.05224428918547
$block2[5,5] = 1.015812029323789
$block2[6,0] = 1.037116965134417
$block2[6,1] = 1.034845224631493
$block2[6,2] = 1.040752357939634
$block2[6,3] = 1.041765441272747
$block2[6,4] = 1.051393430953514
$block2[6,5] = 1.015613989462467
$block2[7,0] = 1.036775932333466
$block2[7,1] = 1.033802321611222
$block2[7,2] = 1.040041481260421
$block2[7,3] = 1.04051869612335
$block2[7,4] = 1.049898480186861
$block2[7,5] = 1.015264887336929
$block2[8,0] = 1.036543940521159
$block2[8,1] = 1.03310747841982
$block2[8,2] = 1.039565607553862
$block2[8,3] = 1.039689626695564
$block2[8,4] = 1.048904888293802
$block2[8,5] = 1.01503209321463
$block2[9,0] = 1.036442395647424
$block2[9,1] = 1.032806720879955
$block2[9,2] = 1.039359102515786
$block2[9,3] = 1.039331146091661
$block2[9,4] = 1.048475395505545
$block2[9,5] = 1.014931282720544
$block2[10,0] = 1.036404514254925
$block2[10,1] = 1.03269502460549
$block2[10,2] = 1.039282331368927
$block2[10,3] = 1.03919806888158
$block2[10,4] = 1.048315975578092
$block2[10,5] = 1.01489383633266
$block2[11,0] = 1.03641264731039
$block2[11,1] = 1.032718982978693
$block2[11,2] = 1.039298802003655
$block2[11,3] = 1.03922661082031
$block2[11,4] = 1.048350166559863
$block2[11,5] = 1.014901868742084
$block2[12,0] = 1.036439267682434
$block2[12,1] = 1.032797487646721
$block2[12,2] = 1.039352757926669
$block2[12,3] = 1.039320144278775
$block2[12,4] = 1.048462215492867
$block2[12,5] = 1.01492818740562
$block2[13,0] = 1.036455647782824
$block2[13,1] = 1.032845859450155
$block2[13,2] = 1.039385993250699
$block2[13,3] = 1.039377783770524
$block2[13,4] = 1.048531267541258
$block2[13,5] = 1.014944403100529
$block2[14,0] = 1.036550656803527
$block2[14,1] = 1.033127441224652
$block2[14,2] = 1.039579303282723
$block2[14,3] = 1.039713428789501
$block2[14,4] = 1.048933408038889
$block2[14,5] = 1.015038783524801
$block2[15,0] = 1.036609961946749
$block2[15,1] = 1.033304101523175
$block2[15,2] = 1.039700442458959
$block2[15,3] = 1.039924108184735
$block2[15,4] = 1.049185859273589
$block2[15,5] = 1.015097983777421
$block2[16,0] = 1.03664444832562
$block2[16,1] = 1.033407155411628
$block2[16,2] = 1.039771057514574
$block2[16,3] = 1.040047043202339
$block2[16,4] = 1.049333180996532
$block2[16,5] = 1.015132513372326
$block2[17,0] = 1.0366561894036
$block2[17,1] = 1.033442295937453
$block2[17,2] = 1.039795128023255
$block2[17,3] = 1.040088969180253
$block2[17,4] = 1.049383425917572
$block2[17,5] = 1.015144286902581
$block2[18,0] = 1.036603609952938
$block2[18,1] = 1.033285146405494
$block2[18,2] = 1.039687449848228
$block2[18,3] = 1.039901499170024
$block2[18,4] = 1.04915876626903
$block2[18,5] = 1.015091632241291
$block2[19,0] = 1.036431433142526
$block2[19,1] = 1.0327743694583
$block2[19,2] = 1.039336871062222
$block2[19,3] = 1.039292598828648
$block2[19,4] = 1.048429216746223
$block2[19,5] = 1.014920437235626
$block2[20,0] = 1.036322235463167
$block2[20,1] = 1.032453331216006
$block2[20,2] = 1.039116067216297
$block2[20,3] = 1.038910213650145
$block2[20,4] = 1.047971173613505
$block2[20,5] = 1.014812795279252
$block2[21,0] = 1.036380212368099
$block2[21,1] = 1.032623509032985
$block2[21,2] = 1.039233155160612
$block2[21,3] = 1.039112879661506
$block2[21,4] = 1.048213928395086
$block2[21,5] = 1.014869858615691
$block2[22,0] = 1.03660648047101
$block2[22,1] = 1.033293711373674
$block2[22,2] = 1.039693320783973
$block2[22,3] = 1.039911715058479
$block2[22,4] = 1.049171008211857
$block2[22,5] = 1.015094502229947
$block2[23,0] = 1.036864918144286
$block2[23,1] = 1.03407186769265
$block2[23,2] = 1.040225610680041
$block2[23,3] = 1.04084064709585
$block2[23,4] = 1.050284432222883
$block2[23,5] = 1.015355151149125
$ws.Range("I2:N25").Value = $block2
